$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12; this shifts the existing rows 12..172
# down to 13..173 (which is exactly the row-by-row downward shift seen in
# the diff) and grows the sheet dimension to A1:R173 automatically.
$ws.Rows.Item(12).Insert()

# Populate the freshly inserted row 12 with the new weekly record.
$ws.Range("A12").Value = 8
$ws.Range("B12").Value = "Terminal La Palmera de La Serena"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 44503
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 100112012
$ws.Range("G12").Value = "Espinaca"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 2700
$ws.Range("K12").Value = 450
$ws.Range("L12").Value = 500
$ws.Range("M12").Value = 475
$ws.Range("N12").Value = '$/atado 300 a 500 gramos'
$ws.Range("O12").Value = 'Provincia del Elquí'
$ws.Range("P12").Value = 950
$ws.Range("Q12").Value = 0.5
$ws.Range("R12").Value = "Hortaliza"
